$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B52").Value = "All columns"
$ws.Range("C52").Value = "Accuracy (A3)"
$ws.Range("D52").Value = "2024-11-29 16:16:51"
$ws.Range("E52").Value = "no threshold"
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = "EwertM"

$ws.Range("A53").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B53").Value = "POPULATION, ESTIMATE_CLASSIFICATION, ESTIMATE_METHOD"
$ws.Range("C53").Value = "Consistency (C1)"
$ws.Range("D53").Value = "2024-11-29 16:17:07"
$ws.Range("E53").Value = 0.91
$ws.Range("F53").Value = 0.9946666666666667
$ws.Range("G53").Value = "EwertM"

$ws.Range("A54").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B54").Value = "POPULATION, ESTIMATE_CLASSIFICATION, ESTIMATE_METHOD"
$ws.Range("C54").Value = "Consistency (C1)"
$ws.Range("D54").Value = "2024-11-29 16:17:18"
$ws.Range("E54").Value = 0.91
$ws.Range("F54").Value = 0.9946666666666667
$ws.Range("G54").Value = "EwertM"

$ws.Range("A55").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B55").Value = "AREA, ANALYSIS_YR, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, NO_INSPECTIONS_USED, ACT_ID, POP_ID, GFE_ID"
$ws.Range("C55").Value = "Accuracy (A1)"
$ws.Range("D55").Value = "2024-11-29 16:17:19"
$ws.Range("E55").Value = "no threshold"
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = "EwertM"

$ws.Range("A56").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B56").Value = "NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER"
$ws.Range("C56").Value = "Accuracy (A2)"
$ws.Range("D56").Value = "2024-11-29 16:17:29"
$ws.Range("E56").Value = 1.5
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = "EwertM"

$ws.Range("A57").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B57").Value = "All columns"
$ws.Range("C57").Value = "Accuracy (A3)"
$ws.Range("D57").Value = "2024-11-29 16:17:40"
$ws.Range("E57").Value = "no threshold"
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = "EwertM"

$ws.Range("A58").Value = "North and Central Coast NuSEDS_20241004"
$ws.Range("B58").Value = "All columns"
$ws.Range("C58").Value = "Accuracy (A3)"
$ws.Range("D58").Value = "2024-11-29 16:18:57"
$ws.Range("E58").Value = "no threshold"
$ws.Range("F58").Value = 1
$ws.Range("G58").Value = "EwertM"

$ws.Range("A59").Value = "West Coast Vancouver Island NuSEDS_20241004"
$ws.Range("B59").Value = "All columns"
$ws.Range("C59").Value = "Accuracy (A3)"
$ws.Range("D59").Value = "2024-11-29 16:19:26"
$ws.Range("E59").Value = "no threshold"
$ws.Range("F59").Value = 1
$ws.Range("G59").Value = "EwertM"

$ws.Range("A60").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B60").Value = "All columns"
$ws.Range("C60").Value = "Accuracy (A3)"
$ws.Range("D60").Value = "2024-11-29 16:20:13"
$ws.Range("E60").Value = "no threshold"
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = "EwertM"

$ws.Range("A61").Value = "Yukon and Transboundary NuSEDS_20241004"
$ws.Range("B61").Value = "All columns"
$ws.Range("C61").Value = "Completeness (C1)"
$ws.Range("D61").Value = "2024-11-29 16:25:08"
$ws.Range("E61").Value = 0.75
$ws.Range("F61").Value = 0.9153491436100132
$ws.Range("G61").Value = "EwertM"

$ws.Range("A62").Value = "West Coast Vancouver Island NuSEDS_20241004"
$ws.Range("B62").Value = "All columns"
$ws.Range("C62").Value = "Completeness (C1)"
$ws.Range("D62").Value = "2024-11-29 16:26:35"
$ws.Range("E62").Value = 0.75
$ws.Range("F62").Value = 0.888227784909428
$ws.Range("G62").Value = "EwertM"

$ws.Range("A63").Value = "North and Central Coast NuSEDS_20241004"
$ws.Range("B63").Value = "All columns"
$ws.Range("C63").Value = "Completeness (C1)"
$ws.Range("D63").Value = "2024-11-29 16:27:50"
$ws.Range("E63").Value = 0.75
$ws.Range("F63").Value = 0.8610714086106374
$ws.Range("G63").Value = "EwertM"

$ws.Range("A64").Value = "Johnstone Strait and Strait of Georgia NuSEDS_20241004"
$ws.Range("B64").Value = "All columns"
$ws.Range("C64").Value = "Completeness (C1)"
$ws.Range("D64").Value = "2024-11-29 16:28:36"
$ws.Range("E64").Value = 0.75
$ws.Range("F64").Value = 0.8482207305966877
$ws.Range("G64").Value = "EwertM"
